# adding info MICS and LSMS
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add remark in column F for row 19 (Year 2007): GIS data not available
$ws.Range("F19").Value = "GIS data not available"

# Add remark in column F for row 23 (Year 2011): GIS data not available
$ws.Range("F23").Value = "GIS data not available"

# Row 28 (Year 2016) -> becomes a text range "2016-2017"
$ws.Range("A28").Value = "2016-2017"

# Add remark in column F for row 33 (Year 2021): Anthropometric data not available
$ws.Range("F33").Value = "Anthropometric data not available"

# Row 22 (Year 2010) -> becomes a text range "2010 - 2011"
$ws.Range("A22").Value = "2010 - 2011"

# Update the selected cell/range shown when the workbook is opened
$ws.Range("L30").Select()
